# Update odds values in the "Jogos da Semana" worksheet (row 6 and row 7)
# as per the latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("H6").Value  = 3.8
$ws.Range("I6").Value  = 8.75
$ws.Range("J6").Value  = 1.93
$ws.Range("K6").Value  = 2.15
$ws.Range("L6").Value  = 7.5
$ws.Range("N6").Value  = 6.8
$ws.Range("Y6").Value  = 8.5
$ws.Range("Z6").Value  = 8.75
$ws.Range("AC6").Value = 8.5
$ws.Range("AE6").Value = 22
$ws.Range("AF6").Value = 120
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 27
$ws.Range("AL6").Value = 120
$ws.Range("AM6").Value = 110
$ws.Range("AO6").Value = 6.5
$ws.Range("AP6").Value = 18.5
$ws.Range("AQ6").Value = 20
$ws.Range("AR6").Value = 55
$ws.Range("AS6").Value = 300
$ws.Range("AT6").Value = 2.47
$ws.Range("AU6").Value = 8.5
$ws.Range("AW6").Value = 9
$ws.Range("AX6").Value = 55
$ws.Range("AY6").Value = 50
$ws.Range("AZ6").Value = 450

# Row 7 updates
$ws.Range("N7").Value = 7.9
